# Weekly refresh of the Fruta/Hortaliza (Guayaba) data:
# the per-day records (Fecha + Volumen + Precio min/max/promedio/$Kg)
# get reshuffled onto different rows as new daily rows roll in/out.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for each data row (row number => Fecha, Volumen, Precio minimo,
# Precio maximo, Precio promedio ponderado, Precio $/Kg). Columns N, O and P
# always mirror S (the $/Kg rate) for this sheet.
$rows = @(
    @{ Row = 2;  Fecha = 44432; Volumen = 30;  Precio = 1300 },
    @{ Row = 3;  Fecha = 44424; Volumen = 50;  Precio = 1200 },
    @{ Row = 4;  Fecha = 44418; Volumen = 40;  Precio = 1200 },
    @{ Row = 5;  Fecha = 44357; Volumen = 35;  Precio = 1000 },
    @{ Row = 6;  Fecha = 44343; Volumen = 60;  Precio = 1300 },
    @{ Row = 7;  Fecha = 44473; Volumen = 120; Precio = 1200 },
    @{ Row = 8;  Fecha = 44435; Volumen = 130; Precio = 1300 },
    @{ Row = 9;  Fecha = 44405; Volumen = 50;  Precio = 1200 },
    @{ Row = 10; Fecha = 44417; Volumen = 80;  Precio = 1200 },
    @{ Row = 11; Fecha = 44476; Volumen = 80;  Precio = 1200 },
    @{ Row = 12; Fecha = 44431; Volumen = 100; Precio = 1300 },
    @{ Row = 13; Fecha = 44438; Volumen = 60;  Precio = 1200 }
)

foreach ($r in $rows) {
    $i = $r.Row
    $ws.Cells.Item($i, 4).Value  = $r.Fecha    # D: Fecha
    $ws.Cells.Item($i, 13).Value = $r.Volumen  # M: Volumen
    $ws.Cells.Item($i, 14).Value = $r.Precio   # N: Precio minimo
    $ws.Cells.Item($i, 15).Value = $r.Precio   # O: Precio maximo
    $ws.Cells.Item($i, 16).Value = $r.Precio   # P: Precio promedio ponderado
    $ws.Cells.Item($i, 19).Value = $r.Precio   # S: Precio $/Kg
}
